# Slide 4 ("Learnings") - Content Placeholder 2:
#  - add two new bullets at the top ("Initially had 3 different..." and
#    "Taking time to understand results" moved up)
#  - add three new lvl=1 sub-bullets (Isa / Reto / Andreas) after the
#    "Team collaboration..." bullet
#  - "Make code OS agnostic (file paths)" ends up as the last bullet
#  - body autofit switches to "shrink text on overflow" (normAutofit)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tf = $sh.TextFrame
$tr = $tf.TextRange

$newText = @(
    "Initially had 3 different approaches, as team decided to go with decision tree => not a black box",
    "Taking time to understand results",
    "Decision tree depth is very important (best results @3 for this case)",
    "Team collaboration is much easier when using a common repository",
    "Isa focused on the graphics and analysis part",
    "Reto focused on the data cleansing",
    "Andreas focused on the code clean up and os independent parts",
    "Make code OS agnostic (file paths)"
) -join "`r"

$tr.Text = $newText

# Demote the three "who did what" bullets to the second outline level.
$tr.Paragraphs(5,1).IndentLevel = 2
$tr.Paragraphs(6,1).IndentLevel = 2
$tr.Paragraphs(7,1).IndentLevel = 2

# "Reto focused on the data cleansing" -> split into two runs: "Reto" / " focused on the data cleansing"
$retoPara = $tr.Paragraphs(6,1)
$retoPara.Text = "Reto"
[void]$retoPara.InsertAfter(" focused on the data cleansing")

# "Andreas focused on the code clean up and os independent parts" -> split into three runs
$andreasPara = $tr.Paragraphs(7,1)
$andreasPara.Text = "Andreas focused on the code clean up and "
[void]$andreasPara.InsertAfter("os")
[void]$andreasPara.InsertAfter(" independent parts")

# Shrink text on overflow (<a:normAutofit/>)
$tf.AutoSize = 2
